$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that need the bold red "highlight" style applied (matches existing style index 2
# used elsewhere in the sheet: bold, red font, default fill/border).
$cellsToStyle = @("A3", "A4", "A11", "D11", "A12", "D12", "B13", "B14")

# Update values first
$ws.Range("A4").Value = "GNZG"
$ws.Range("D11").Value = "STAN"
$ws.Range("A12").Value = "STAN"
$ws.Range("D12").Value = "STAN"
$ws.Range("B14").Value = "ALAST"

# Apply the existing bold-red highlight format (already used on A2/A10/B12/D10)
# to the target cells by copying its format, rather than building new font
# objects (which would otherwise create a duplicate style entry).
$src = $ws.Range("A2")
$src.Copy()
foreach ($addr in $cellsToStyle) {
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
}
